# Reduce the amount of "texture opened" values in columns N:P (rows 2-41)
# by 29, widen those columns to fit the new (shorter) values, and leave the
# selection on N7 (the cell the author was last looking at when they saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 41; $r++) {
    foreach ($col in @("N", "O", "P")) {
        $cell = $ws.Range("$col$r")
        $cell.Value = $cell.Value() - 29
    }
}

# Column widths for N:P no longer use Excel's auto "best fit"; set an
# explicit width instead (matches column M's width).
$ws.Range("N1:P1").EntireColumn.ColumnWidth = $ws.Range("M1").EntireColumn.ColumnWidth

# Restore the last active selection.
$ws.Range("N7").Select()
